$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update column F (dSF) values for the listed rows to match repulled data
$ws.Range("F2").Value  = 2
$ws.Range("F3").Value  = 1
$ws.Range("F4").Value  = -2
$ws.Range("F5").Value  = -5
$ws.Range("F6").Value  = 4
$ws.Range("F7").Value  = -5
$ws.Range("F8").Value  = -4
$ws.Range("F9").Value  = -4
$ws.Range("F11").Value = 1
